$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=""67.077.66"""
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Formula = "=""3.503.09"""
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Formula = "=""594.76"""
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Formula = "=""173.16"""
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Formula = "=""0.600"""
$ws.Range("E8").Value = "  +1.84%  "
$ws.Range("E9").Value = "  +3.94%  "
$ws.Range("D10").Formula = "=""7.27"""
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Formula = "=""0.434"""
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Formula = "=""4.109.03"""
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Formula = "=""29.01"""
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Formula = "=""67.069.58"""
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Formula = "=""3.492.54"""
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Formula = "=""394.34"""
$ws.Range("D21").Formula = "=""8.01"""
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Formula = "=""73.10"""
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Formula = "=""0.537"""
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Formula = "=""5.70"""
$ws.Range("E25").Value = "  -3.00%  "
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("D27").Formula = "=""10.28"""
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Formula = "=""0.997"""
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Formula = "=""23.76"""
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").Formula = "=""163.51"""
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Formula = "=""6.98"""
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Formula = "=""0.0745"""
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Formula = "=""27.37"""
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").Formula = "=""26.36"""
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Formula = "=""2.811.84"""
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("D46").Formula = "=""42.68"""
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").Formula = "=""335.94"""
$ws.Range("E48").Value = "  -4.57%  "
$ws.Range("D49").Formula = "=""34.44"""
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("E51").Value = "  +0.31%  "

# Convert the text-literal formulas above into plain static text values
# (so the cells keep their original "text" representation/style instead of a formula).
$used = $ws.UsedRange
$used.Copy()
$used.PasteSpecial(-4163)
$excel.CutCopyMode = 0

